$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 2: Boyd Gaming (BYD) ---
$ws.Range("F2").Value = 6285.69
$ws.Range("I2").Value = 65.66
$ws.Range("J2").Value = 6.21
$ws.Range("K2").Value = 6.09
$ws.Range("L2").Value = 6.6

# --- Row 3: Churchill Downs (CHDN) ---
$ws.Range("F3").Value = 8806.6
$ws.Range("I3").Value = 119.51
$ws.Range("J3").Value = 5.44
$ws.Range("K3").Value = 5.2

# --- Row 4: DraftKings (DKNG) ---
$ws.Range("F4").Value = 39319.16
$ws.Range("I4").Value = 45.37
$ws.Range("J4").Value = -0.24
$ws.Range("K4").Value = -1.54
# M4's EG1 formula (K4/J4-1) is removed/cleared in this edit
$ws.Range("M4").ClearContents()

# --- Row 5: International Game Technology (IGT) ---
$ws.Range("F5").Value = 4222.46
$ws.Range("I5").Value = 21.21
$ws.Range("J5").Value = 1.48
$ws.Range("K5").Value = 1.2
$ws.Range("L5").Value = 2.19

# --- Row 6: Light & Wonder (LNW) ---
$ws.Range("F6").Value = 8923.64
$ws.Range("I6").Value = 99.34
$ws.Range("J6").Value = 3.78
$ws.Range("K6").Value = 1.79
$ws.Range("L6").Value = 5.11

# --- Row 7: Las Vegas Sands (LVS) ---
$ws.Range("F7").Value = 39509.07
$ws.Range("I7").Value = 52.59
$ws.Range("J7").Value = 2.74
$ws.Range("K7").Value = 1.91

# --- Row 8: MGM Resorts International (MGM) ---
$ws.Range("F8").Value = 14538.34
$ws.Range("I8").Value = 45.86
$ws.Range("J8").Value = 2.48
$ws.Range("K8").Value = 2.28

# --- Row 9: Melco Resorts & Entertainment Limited (MLCO) ---
$ws.Range("F9").Value = 3213.54
$ws.Range("I9").Value = 7.22
$ws.Range("J9").Value = 0.28
$ws.Range("K9").Value = -0.19
$ws.Range("L9").Value = 0.51
# EG1/EG2 lose their formulas and become plain values
$ws.Range("M9").Value = -1
$ws.Range("N9").Value = 1

# --- Row 10: Roblox (RBLX) ---
$ws.Range("F10").Value = 21401.11
$ws.Range("I10").Value = 36.8
$ws.Range("J10").Value = -2.1
$ws.Range("K10").Value = -1.93

# --- Row 11: Red Rock Resorts (RRR) ---
$ws.Range("F11").Value = 6295.51
$ws.Range("I11").Value = 59.85
$ws.Range("J11").Value = 1.95
$ws.Range("K11").Value = 1.57
$ws.Range("L11").Value = 2.19

# --- Row 12: Wynn Resorts (WYNN) ---
$ws.Range("F12").Value = 11893.75
$ws.Range("I12").Value = 106.12
$ws.Range("J12").Value = 5.28
$ws.Range("K12").Value = 3.24
$ws.Range("L12").Value = 5.67

# --- Selection moved from T3 to H20 ---
$null = $ws.Range("H20").Select()
